$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $t = $cell.Text
    if ($t -eq $oldValue) {
        $cell.Value = $newValue
    }
}
